# Fixed Asset test-data workbook: scrub the committed Oracle Cloud login
# credentials that were left on the "Input_Value" sheet (URL / UserName /
# Password, stored in J2:L2) before re-uploading the file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input_Value")

# Remove the stored URL / UserName / Password values.
$ws.Range("J2:L2").ClearContents()

# Leave the view scrolled over to, and selecting, the cells that were edited
# (mirrors the saved sheet view: columns scrolled right, J2:L2 selected).
$ws.Activate()
$ws.Range("J2:L2").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 5
